$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (dates as Excel serial numbers, matching the
# existing column A date format/style already used by the sheet).
$newRows = @(
    @{ Row = 234; Date = 44308; B = 0; C = 11; D = 294.4325481798715 },
    @{ Row = 235; Date = 44309; B = 0; C = 9;  D = 240.8993576017131 },
    @{ Row = 236; Date = 44310; B = 5; C = 13; D = 347.9657387580299 },
    @{ Row = 237; Date = 44311; B = 0; C = 10; D = 267.6659528907923 },
    @{ Row = 238; Date = 44312; B = 0; C = 10; D = 267.6659528907923 }
)

$lastExistingRow = 233

foreach ($item in $newRows) {
    $r = $item.Row

    # Copy the style (border/font/alignment/number-format) from the cell
    # directly above so the new date cell matches the rest of column A.
    $ws.Range("A$lastExistingRow").Copy($ws.Range("A$r"))

    $ws.Range("A$r").Value = $item.Date
    $ws.Range("B$r").Value = $item.B
    $ws.Range("C$r").Value = $item.C
    $ws.Range("D$r").Value = $item.D

    $lastExistingRow = $r
}
